# "Builder and cases adapted with battery block"
#
# This script:
#   1. Updates SolarPV (PV_forecast) values.
#   2. Updates Source (Irrigation1_Q) values and makes it the active sheet.
#   3. Tweaks the HydroSwitch sheet selection.
#   4. Adds a bold style to Grid!C24 (new formatted cell below the existing data).
#   5. Appends two new worksheets (Battery_Ex0, Battery) with a small
#      Battery_E / Battery_P time-series, at the end of the workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. SolarPV: update the PV_forecast series (A2:A6)
# ---------------------------------------------------------------------
$solarPV = $wb.Worksheets.Item("SolarPV")
$solarPV.Cells.Item(2, 1).Value = 0.4
$solarPV.Cells.Item(3, 1).Value = 1
$solarPV.Cells.Item(4, 1).Value = 1.7
$solarPV.Cells.Item(5, 1).Value = 1.2
$solarPV.Cells.Item(6, 1).Value = 0.5
[void]$solarPV.Range("A7").Select()

# ---------------------------------------------------------------------
# 2. Source: update the Irrigation1_Q series (A2:A6) and select it
# ---------------------------------------------------------------------
$source = $wb.Worksheets.Item("Source")
$source.Cells.Item(3, 1).Value = 0
$source.Cells.Item(4, 1).Value = 0
$source.Cells.Item(5, 1).Value = 0
$source.Cells.Item(6, 1).Value = 2

# Source becomes the active/selected sheet (was SolarPV before)
[void]$source.Activate()
[void]$source.Range("A7").Select()

# ---------------------------------------------------------------------
# 3. HydroSwitch: move the selection
# ---------------------------------------------------------------------
$hydroSwitch = $wb.Worksheets.Item("HydroSwitch")
[void]$hydroSwitch.Range("I13").Select()

# ---------------------------------------------------------------------
# 4. Grid: add a new bold-styled cell at C24
# ---------------------------------------------------------------------
$grid = $wb.Worksheets.Item("Grid")
$grid.Range("C24").Font.Bold = $true

# ---------------------------------------------------------------------
# 5. Append the Battery_Ex0 / Battery worksheets at the end
# ---------------------------------------------------------------------
function Add-BatterySheet {
    param([string]$sheetName, [string]$activeCellRef)

    $lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
    $newSheet = $wb.Worksheets.Add($null, $lastSheet)
    $newSheet.Name = $sheetName

    $newSheet.Range("A1").Value = "Battery_E"
    $newSheet.Range("B1").Value = "Battery_P"

    for ($row = 2; $row -le 6; $row++) {
        $newSheet.Cells.Item($row, 1).Value = 1000
        $newSheet.Cells.Item($row, 2).Value = 1000
    }

    [void]$newSheet.Range($activeCellRef).Select()
}

Add-BatterySheet "Battery_Ex0" "E12"
Add-BatterySheet "Battery" "A2"

# Re-assert Source as the active sheet/tab (adding sheets above moved focus)
[void]$source.Activate()
[void]$source.Range("A7").Select()
